$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 19
$ws_ALC.Range("H19").Value = 267.45947
$ws_ALC.Range("I19").Value = 227.26086
$ws_ALC.Range("J19").Value = 333.5
$ws_ALC.Range("K19").Value = 227.26086
$ws_ALC.Range("L19").Value = 333.5
$ws_ALC.Range("M19").Value = -52.26086000000001
$ws_ALC.Range("N19").Value = -683.5

# ALC row 33
$ws_ALC.Range("H33").Value = 249.14285
$ws_ALC.Range("I33").Value = 252.8
$ws_ALC.Range("J33").Value = 240
$ws_ALC.Range("K33").Value = 252.8
$ws_ALC.Range("L33").Value = 240
$ws_ALC.Range("M33").Value = -23.80000000000001
$ws_ALC.Range("N33").Value = -698

# ALC row 76
$ws_ALC.Range("H76").Value = 3670.8235
$ws_ALC.Range("I76").Value = 3675.25
$ws_ALC.Range("K76").Value = 3675.25
$ws_ALC.Range("M76").Value = -3360.25

# ALC row 79
$ws_ALC.Range("H79").Value = 3670.8235
$ws_ALC.Range("I79").Value = 3675.25
$ws_ALC.Range("K79").Value = 3675.25
$ws_ALC.Range("M79").Value = -2583.25

# ALC row 132
$ws_ALC.Range("H132").Value = 1230.1613
$ws_ALC.Range("I132").Value = 1246.0344
$ws_ALC.Range("J132").Value = 1000
$ws_ALC.Range("K132").Value = 3738.1032
$ws_ALC.Range("L132").Value = 3000
$ws_ALC.Range("M132").Value = -1208.1032
$ws_ALC.Range("N132").Value = -8060

# ALC row 137
$ws_ALC.Range("H137").Value = 5035
$ws_ALC.Range("I137").Value = 4870.0586
$ws_ALC.Range("J137").Value = 5199.9414
$ws_ALC.Range("K137").Value = 14610.1758
$ws_ALC.Range("L137").Value = 15599.8242
$ws_ALC.Range("M137").Value = -12060.1758
$ws_ALC.Range("N137").Value = -20699.8242

# ARM row 32
$ws_ARM.Range("H32").Value = 2796.12
$ws_ARM.Range("I32").Value = 2468.1099
$ws_ARM.Range("J32").Value = 6112.6665
$ws_ARM.Range("K32").Value = 2468.1099
$ws_ARM.Range("L32").Value = 6112.6665
$ws_ARM.Range("M32").Value = -2181.1099
$ws_ARM.Range("N32").Value = -6686.6665

# ARM row 63
$ws_ARM.Range("H63").Value = 2696.147
$ws_ARM.Range("I63").Value = 2519.9583
$ws_ARM.Range("J63").Value = 3119
$ws_ARM.Range("K63").Value = 2519.9583
$ws_ARM.Range("L63").Value = 3119
$ws_ARM.Range("M63").Value = -1833.9583
$ws_ARM.Range("N63").Value = -4491

# ARM row 66
$ws_ARM.Range("H66").Value = 2696.147
$ws_ARM.Range("I66").Value = 2519.9583
$ws_ARM.Range("J66").Value = 3119
$ws_ARM.Range("K66").Value = 12599.7915
$ws_ARM.Range("L66").Value = 15595
$ws_ARM.Range("M66").Value = -9167.791499999999
$ws_ARM.Range("N66").Value = -22459

# ARM row 88
$ws_ARM.Range("H88").Value = 2639.1428
$ws_ARM.Range("J88").Value = 2607
$ws_ARM.Range("L88").Value = 2607
$ws_ARM.Range("N88").Value = -3419

# ARM row 91
$ws_ARM.Range("H91").Value = 2639.1428
$ws_ARM.Range("J91").Value = 2607
$ws_ARM.Range("L91").Value = 2607
$ws_ARM.Range("N91").Value = -5415

# ARM row 97
$ws_ARM.Range("H97").Value = 1933.3334
$ws_ARM.Range("I97").Value = 1781.5834
$ws_ARM.Range("J97").Value = 2236.8333
$ws_ARM.Range("K97").Value = 1781.5834
$ws_ARM.Range("L97").Value = 2236.8333
$ws_ARM.Range("M97").Value = -1285.5834
$ws_ARM.Range("N97").Value = -3228.8333

# BSM row 29
$ws_BSM.Range("H29").Value = 0
$ws_BSM.Range("I29").Value = 0
$ws_BSM.Range("K29").Value = 0
$ws_BSM.Range("M29").ClearContents()

# BSM row 99
$ws_BSM.Range("H99").Value = 1136.5
$ws_BSM.Range("I99").Value = 1023.31036
$ws_BSM.Range("J99").Value = 1501.2222
$ws_BSM.Range("K99").Value = 1023.31036
$ws_BSM.Range("L99").Value = 1501.2222
$ws_BSM.Range("M99").Value = 474.6896400000001
$ws_BSM.Range("N99").Value = -4497.2222

# BSM row 105
$ws_BSM.Range("H105").Value = 3506.6667
$ws_BSM.Range("I105").Value = 3506.6667
$ws_BSM.Range("J105").Value = 0
$ws_BSM.Range("K105").Value = 3506.6667
$ws_BSM.Range("L105").Value = 0
$ws_BSM.Range("M105").Value = -1759.6667
$ws_BSM.Range("N105").ClearContents()

# BSM row 134
$ws_BSM.Range("H134").Value = 1116.1818
$ws_BSM.Range("I134").Value = 993.7931
$ws_BSM.Range("J134").Value = 2003.5
$ws_BSM.Range("K134").Value = 2981.3793
$ws_BSM.Range("L134").Value = 6010.5
$ws_BSM.Range("M134").Value = -446.3793000000001
$ws_BSM.Range("N134").Value = -11080.5

# CRP row 31
$ws_CRP.Range("H31").Value = 3789.8696
$ws_CRP.Range("I31").Value = 5602.4
$ws_CRP.Range("J31").Value = 3286.389
$ws_CRP.Range("K31").Value = 5602.4
$ws_CRP.Range("L31").Value = 3286.389
$ws_CRP.Range("M31").Value = -5307.4
$ws_CRP.Range("N31").Value = -3876.389

# CRP row 34
$ws_CRP.Range("H34").Value = 3789.8696
$ws_CRP.Range("I34").Value = 5602.4
$ws_CRP.Range("J34").Value = 3286.389
$ws_CRP.Range("K34").Value = 5602.4
$ws_CRP.Range("L34").Value = 3286.389
$ws_CRP.Range("M34").Value = -5400.4
$ws_CRP.Range("N34").Value = -3690.389

# CRP row 132
$ws_CRP.Range("H132").Value = 1753.4688
$ws_CRP.Range("I132").Value = 1363.591
$ws_CRP.Range("K132").Value = 4090.773
$ws_CRP.Range("M132").Value = -1560.773

# CRP row 134
$ws_CRP.Range("H134").Value = 2224.15
$ws_CRP.Range("I134").Value = 2398.0625
$ws_CRP.Range("J134").Value = 1528.5
$ws_CRP.Range("K134").Value = 7194.1875
$ws_CRP.Range("L134").Value = 4585.5
$ws_CRP.Range("M134").Value = -4659.1875
$ws_CRP.Range("N134").Value = -9655.5

# CUL row 113
$ws_CUL.Range("H113").Value = 1115.591
$ws_CUL.Range("I113").Value = 1889.4
$ws_CUL.Range("J113").Value = 470.75
$ws_CUL.Range("K113").Value = 5668.200000000001
$ws_CUL.Range("L113").Value = 1412.25
$ws_CUL.Range("M113").Value = -3498.200000000001
$ws_CUL.Range("N113").Value = -5752.25

# GSM row 70
$ws_GSM.Range("H70").Value = 4137.2354
$ws_GSM.Range("I70").Value = 4126.875
$ws_GSM.Range("J70").Value = 4146.4443
$ws_GSM.Range("K70").Value = 4126.875
$ws_GSM.Range("L70").Value = 4146.4443
$ws_GSM.Range("M70").Value = -3856.875
$ws_GSM.Range("N70").Value = -4686.4443

# GSM row 73
$ws_GSM.Range("H73").Value = 4137.2354
$ws_GSM.Range("I73").Value = 4126.875
$ws_GSM.Range("J73").Value = 4146.4443
$ws_GSM.Range("K73").Value = 4126.875
$ws_GSM.Range("L73").Value = 4146.4443
$ws_GSM.Range("M73").Value = -3190.875
$ws_GSM.Range("N73").Value = -6018.4443

# GSM row 122
$ws_GSM.Range("H122").Value = 2278.4644
$ws_GSM.Range("I122").Value = 1372.238
$ws_GSM.Range("J122").Value = 4997.143
$ws_GSM.Range("K122").Value = 4116.714
$ws_GSM.Range("L122").Value = 14991.429
$ws_GSM.Range("M122").Value = -1666.714
$ws_GSM.Range("N122").Value = -19891.429

# LTW row 7
$ws_LTW.Range("H7").Value = 2131
$ws_LTW.Range("I7").Value = 1822.75
$ws_LTW.Range("K7").Value = 1822.75
$ws_LTW.Range("M7").Value = -1710.75

# LTW row 46
$ws_LTW.Range("H46").Value = 261935.6
$ws_LTW.Range("I46").Value = 1184.8182
$ws_LTW.Range("J46").Value = 979000.25
$ws_LTW.Range("K46").Value = 1184.8182
$ws_LTW.Range("L46").Value = 979000.25
$ws_LTW.Range("M46").Value = -996.8181999999999
$ws_LTW.Range("N46").Value = -979376.25

# LTW row 126
$ws_LTW.Range("H126").Value = 2131
$ws_LTW.Range("I126").Value = 1822.75
$ws_LTW.Range("K126").Value = 5468.25
$ws_LTW.Range("M126").Value = -2998.25

# LTW row 132
$ws_LTW.Range("H132").Value = 6207.7646
$ws_LTW.Range("I132").Value = 3474.1428
$ws_LTW.Range("J132").Value = 10623.615
$ws_LTW.Range("K132").Value = 10422.4284
$ws_LTW.Range("L132").Value = 31870.845
$ws_LTW.Range("M132").Value = -7892.428400000001
$ws_LTW.Range("N132").Value = -36930.845

# LTW row 136
$ws_LTW.Range("H136").Value = 2661.1072
$ws_LTW.Range("I136").Value = 2000.6111
$ws_LTW.Range("J136").Value = 3850
$ws_LTW.Range("K136").Value = 6001.8333
$ws_LTW.Range("L136").Value = 11550
$ws_LTW.Range("M136").Value = -3451.8333
$ws_LTW.Range("N136").Value = -16650

# WVR row 24
$ws_WVR.Range("H24").Value = 5000
$ws_WVR.Range("J24").Value = 5000
$ws_WVR.Range("L24").Value = 5000
$ws_WVR.Range("N24").Value = -5460
